$d = $word.ActiveDocument

# The last paragraph in the document is an empty "ListParagraph" styled
# paragraph that holds the "_GoBack" bookmark. It needs to be split into
# two empty paragraphs:
#   1) the original paragraph (keeps the bookmark) loses the ListParagraph
#      style / the 1080-twip left indent, reverting to plain body text
#      formatting (its paragraph-mark run formatting - fonts/size - is left
#      untouched).
#   2) a brand new empty paragraph appended right after it, which keeps the
#      ListParagraph style and the 1080-twip (54pt) left indent that the
#      original paragraph used to have.

$last = $d.Paragraphs.Last
$lastRange = $last.Range

# Step 1: duplicate the paragraph (with its current ListParagraph style +
# indent + run formatting) by inserting a new paragraph mark right after it.
# The new paragraph automatically inherits the current formatting, which is
# exactly what the second (new) paragraph in the target needs.
$lastRange.InsertParagraphAfter()

# Step 2: the original paragraph (still holding the bookmark) is now back
# at the same index it always was; strip its list style / indent.
$first = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$firstRange = $first.Range
$firstRange.Select()
$word.Selection.Style = "Normal"
$word.Selection.ParagraphFormat.LeftIndent = 0
